# Verify_24V_Load_On_Changing_CPU_Of_Panel.xlsx
#
# Updates the "Add Panels" sheet:
#  - Expected 24V PSU Load value (H9) and Expected 2nd 24V PSU Load value
#    (J9) change from "0.000" to "0.205" (text cells, quote-prefixed so
#    they keep their text formatting rather than becoming numbers).
#  - The saved view moves off the "topLeftCell=C1 / J9 selected" state to
#    a 90% zoom with C13 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# --- Data update -----------------------------------------------------
# Leading apostrophe keeps these as text (t="s", quotePrefix style),
# matching the existing "0.000" text values rather than converting the
# cells to numeric.
$ws.Range("H9").Value = "'0.205"
$ws.Range("J9").Value = "'0.205"

# --- View update -------------------------------------------------------
$ws.Activate()
$ws.Range("C13").Select()
$excel.ActiveWindow.Zoom = 90
